$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.44
$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 5.6
$ws.Range("I2").Value = 8.199999999999999
$ws.Range("J2").Value = 4.6
$ws.Range("K2").Value = 5.7
